$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "59.594.95"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "2.670.24"
$ws.Range("E3").Value = "  +1.10%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "525.98"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("E6").Value = "  +0.70%  "

Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.32%  "

Set-TextValue $ws.Range("D9") "6.96"
$ws.Range("E9").Value = "  +10.30%  "

$ws.Range("E10").Value = "  -1.55%  "

Set-TextValue $ws.Range("D11") "0.339"
$ws.Range("E11").Value = "  +0.56%  "

Set-TextValue $ws.Range("D12") "0.131"
$ws.Range("E12").Value = "  +2.25%  "

$ws.Range("D13").Value = "3.131.81"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").Value = "59.603.31"
$ws.Range("E14").Value = "  +0.48%  "

Set-TextValue $ws.Range("D15") "21.34"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "2.676.75"
$ws.Range("E17").Value = "  +1.43%  "

Set-TextValue $ws.Range("D18") "343.21"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("E19").Value = "  -1.23%  "

Set-TextValue $ws.Range("D20") "10.50"
$ws.Range("E20").Value = "  +1.83%  "

$ws.Range("E21").Value = "  +2.72%  "

Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  +2.45%  "

# Rows 24/25: Kaspa and Polygon swapped positions
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D24") "0.416"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D25") "0.169"
$ws.Range("E25").Value = "  +2.18%  "

Set-TextValue $ws.Range("D26") "0.996"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +1.17%  "

Set-TextValue $ws.Range("D28") "7.19"
$ws.Range("E28").Value = "  +0.41%  "

Set-TextValue $ws.Range("D29") "6.69"
$ws.Range("E29").Value = "  +2.11%  "

Set-TextValue $ws.Range("D30") "0.998"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  +1.25%  "

Set-TextValue $ws.Range("D32") "18.95"
$ws.Range("E32").Value = "  +0.20%  "

Set-TextValue $ws.Range("D33") "149.24"
$ws.Range("E33").Value = "  -0.31%  "

Set-TextValue $ws.Range("D34") "4.23"
$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("E35").Value = "  +3.36%  "

Set-TextValue $ws.Range("D36") "0.906"
$ws.Range("E36").Value = "  -4.30%  "

Set-TextValue $ws.Range("D37") "0.902"
$ws.Range("E37").Value = "  +4.49%  "

$ws.Range("E38").Value = "  +0.99%  "

Set-TextValue $ws.Range("D39") "36.95"
$ws.Range("E39").Value = "  +0.86%  "

Set-TextValue $ws.Range("D40") "3.63"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("E41").Value = "  +4.49%  "

Set-TextValue $ws.Range("D42") "20.12"
$ws.Range("E42").Value = "  +2.55%  "

Set-TextValue $ws.Range("D43") "0.999"
$ws.Range("E43").Value = "  +0.10%  "

Set-TextValue $ws.Range("D44") "276.58"
$ws.Range("E44").Value = "  -0.58%  "

Set-TextValue $ws.Range("D45") "0.0979"
$ws.Range("E45").Value = "  -1.06%  "

Set-TextValue $ws.Range("D46") "0.0545"
$ws.Range("E46").Value = "  +2.89%  "

Set-TextValue $ws.Range("D47") "4.94"
$ws.Range("E47").Value = "  +4.67%  "

$ws.Range("D48").Value = "2.071.73"
$ws.Range("E48").Value = "  -0.34%  "

Set-TextValue $ws.Range("D49") "10.53"
$ws.Range("E49").Value = "  +2.04%  "

Set-TextValue $ws.Range("D50") "19.30"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("E51").Value = "  -0.13%  "
